$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set cell A2 to the same value as A3 ("x"), reusing the shared string
$ws.Range("A2").Value = $ws.Range("A3").Value2

# Move the active cell selection from C6 to C5
$ws.Range("C5").Select()
